# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 4; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 5; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 8; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 10; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 12; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 15; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 20; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 22; I = 'ba'; J = 'Appreciation' },
    @{ Row = 25; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 37; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 40; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 46; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 54; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 66; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 81; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 89; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 104; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 120; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 124; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 127; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 135; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 136; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 144; I = 'ba'; J = 'Appreciation' },
    @{ Row = 145; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 152; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 159; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 186; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 187; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 188; I = 'ba'; J = 'Appreciation' },
    @{ Row = 193; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 196; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 197; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 199; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 216; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 221; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 233; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 235; I = 'ba'; J = 'Appreciation' },
    @{ Row = 240; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 253; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 273; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 300; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 310; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 315; I = 'sd'; J = 'Statement-non-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
